$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1268.5294
$ws.Range("I19").Value = 714
$ws.Range("J19").Value = 2599.4
$ws.Range("K19").Value = 714
$ws.Range("L19").Value = 2599.4
$ws.Range("M19").Value = -539
$ws.Range("N19").Value = -2949.4

$ws.Range("H40").Value = 4999.8335
$ws.Range("J40").Value = 4999.8335
$ws.Range("L40").Value = 4999.8335
$ws.Range("N40").Value = -5349.8335

$ws.Range("H125").Value = 28201.75
$ws.Range("I125").Value = 81245.5
$ws.Range("K125").Value = 731209.5
$ws.Range("M125").Value = -728749.5

$ws.Range("H132").Value = 1409.5111
$ws.Range("I132").Value = 937.4146
$ws.Range("K132").Value = 2812.2438
$ws.Range("M132").Value = -282.2437999999997

$ws.Range("H135").Value = 964.05
$ws.Range("I135").Value = 757.75
$ws.Range("K135").Value = 6819.75
$ws.Range("M135").Value = -4284.75

$ws.Range("H137").Value = 2321.8865
$ws.Range("I137").Value = 1999.5483
$ws.Range("J137").Value = 3090.5386
$ws.Range("K137").Value = 5998.644899999999
$ws.Range("L137").Value = 9271.6158
$ws.Range("M137").Value = -3448.644899999999
$ws.Range("N137").Value = -14371.6158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3357.8096
$ws.Range("I2").Value = 3472.2727
$ws.Range("K2").Value = 3472.2727
$ws.Range("M2").Value = -3359.2727

$ws.Range("H45").Value = 407077.62
$ws.Range("I45").Value = 596801.4399999999
$ws.Range("J45").Value = 3914.625
$ws.Range("K45").Value = 596801.4399999999
$ws.Range("L45").Value = 3914.625
$ws.Range("M45").Value = -596424.4399999999
$ws.Range("N45").Value = -4668.625

$ws.Range("H61").Value = 1085.7368
$ws.Range("I61").Value = 1101.8379
$ws.Range("K61").Value = 1101.8379
$ws.Range("M61").Value = -889.8379

$ws.Range("H74").Value = 1417.7544
$ws.Range("I74").Value = 1314.6531
$ws.Range("K74").Value = 1314.6531
$ws.Range("M74").Value = -440.6531

$ws.Range("H77").Value = 1417.7544
$ws.Range("I77").Value = 1314.6531
$ws.Range("K77").Value = 6573.2655
$ws.Range("M77").Value = -2205.2655

$ws.Range("H116").Value = 3357.8096
$ws.Range("I116").Value = 3472.2727
$ws.Range("K116").Value = 3472.2727
$ws.Range("M116").Value = -1178.2727

$ws.Range("H136").Value = 1085.7368
$ws.Range("I136").Value = 1101.8379
$ws.Range("K136").Value = 3305.5137
$ws.Range("M136").Value = -755.5137

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3357.8096
$ws.Range("I3").Value = 3472.2727
$ws.Range("K3").Value = 3472.2727
$ws.Range("M3").Value = -3358.2727

$ws.Range("H22").Value = 620
$ws.Range("I22").Value = 620
$ws.Range("K22").Value = 620
$ws.Range("M22").Value = -447

$ws.Range("H86").Value = 2928.6333
$ws.Range("I86").Value = 2271.818
$ws.Range("J86").Value = 3308.8948
$ws.Range("K86").Value = 2271.818
$ws.Range("L86").Value = 3308.8948
$ws.Range("M86").Value = -1148.818
$ws.Range("N86").Value = -5554.8948

$ws.Range("H89").Value = 2928.6333
$ws.Range("I89").Value = 2271.818
$ws.Range("J89").Value = 3308.8948
$ws.Range("K89").Value = 11359.09
$ws.Range("L89").Value = 16544.474
$ws.Range("M89").Value = -5743.09
$ws.Range("N89").Value = -27776.474

$ws.Range("H107").Value = 1380.7097
$ws.Range("I107").Value = 1343.8334
$ws.Range("J107").Value = 1507.1428
$ws.Range("K107").Value = 1343.8334
$ws.Range("L107").Value = 1507.1428
$ws.Range("M107").Value = 576.1666
$ws.Range("N107").Value = -5347.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3194.8845
$ws.Range("I31").Value = 1750.3
$ws.Range("K31").Value = 1750.3
$ws.Range("M31").Value = -1455.3

$ws.Range("H34").Value = 3194.8845
$ws.Range("I34").Value = 1750.3
$ws.Range("K34").Value = 1750.3
$ws.Range("M34").Value = -1548.3

$ws.Range("H68").Value = 25000
$ws.Range("J68").Value = 25000
$ws.Range("L68").Value = 25000
$ws.Range("N68").Value = -26498

$ws.Range("H71").Value = 25000
$ws.Range("J71").Value = 25000
$ws.Range("L71").Value = 75000
$ws.Range("N71").Value = -82488

$ws.Range("H74").Value = 38701.668
$ws.Range("J74").Value = 56663
$ws.Range("L74").Value = 56663
$ws.Range("N74").Value = -58411

$ws.Range("H77").Value = 38701.668
$ws.Range("J77").Value = 56663
$ws.Range("L77").Value = 169989
$ws.Range("N77").Value = -178725

$ws.Range("H99").Value = 9517.925999999999
$ws.Range("I99").Value = 11560.1875
$ws.Range("J99").Value = 6547.364
$ws.Range("K99").Value = 11560.1875
$ws.Range("L99").Value = 6547.364
$ws.Range("M99").Value = -10062.1875
$ws.Range("N99").Value = -9543.364

$ws.Range("H126").Value = 9517.925999999999
$ws.Range("I126").Value = 11560.1875
$ws.Range("J126").Value = 6547.364
$ws.Range("K126").Value = 34680.5625
$ws.Range("L126").Value = 19642.092
$ws.Range("M126").Value = -32210.5625
$ws.Range("N126").Value = -24582.092

$ws.Range("H132").Value = 1613.6046
$ws.Range("J132").Value = 4000
$ws.Range("L132").Value = 12000
$ws.Range("N132").Value = -17060

$ws.Range("H134").Value = 2361.682
$ws.Range("I134").Value = 2283.6667
$ws.Range("J134").Value = 4000
$ws.Range("K134").Value = 6851.000100000001
$ws.Range("L134").Value = 12000
$ws.Range("M134").Value = -4316.000100000001
$ws.Range("N134").Value = -17070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 802.25
$ws.Range("I107").Value = 916.55554
$ws.Range("J107").Value = 655.2857
$ws.Range("K107").Value = 2749.66662
$ws.Range("L107").Value = 1965.8571
$ws.Range("M107").Value = -829.66662
$ws.Range("N107").Value = -5805.8571

$ws.Range("H109").Value = 862.1111
$ws.Range("I109").Value = 862.1111
$ws.Range("K109").Value = 2586.3333
$ws.Range("M109").Value = -1546.3333

$ws.Range("H113").Value = 714
$ws.Range("I113").Value = 499.75
$ws.Range("J113").Value = 999.6667
$ws.Range("K113").Value = 1499.25
$ws.Range("L113").Value = 2999.0001
$ws.Range("M113").Value = 670.75
$ws.Range("N113").Value = -7339.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3387.5518
$ws.Range("I122").Value = 3365.6785
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 10097.0355
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -7647.0355
$ws.Range("N122").Value = -16900

$ws.Range("H132").Value = 2615.682
$ws.Range("I132").Value = 2377.6667
$ws.Range("J132").Value = 3125.7144
$ws.Range("K132").Value = 7133.000100000001
$ws.Range("L132").Value = 9377.143199999999
$ws.Range("M132").Value = -4603.000100000001
$ws.Range("N132").Value = -14437.1432

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1147.25
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -705

$ws.Range("H27").Value = 1147.25
$ws.Range("I27").Value = 1000
$ws.Range("K27").Value = 1000
$ws.Range("M27").Value = -893

$ws.Range("H61").Value = 1343.2778
$ws.Range("I61").Value = 1249.375
$ws.Range("K61").Value = 1249.375
$ws.Range("M61").Value = -1047.375

$ws.Range("H68").Value = 3047.5
$ws.Range("I68").Value = 3047.5
$ws.Range("K68").Value = 3047.5
$ws.Range("M68").Value = -2298.5

$ws.Range("H71").Value = 3047.5
$ws.Range("I71").Value = 3047.5
$ws.Range("K71").Value = 15237.5
$ws.Range("M71").Value = -11493.5

$ws.Range("H82").Value = 1189
$ws.Range("J82").Value = 999.3333
$ws.Range("L82").Value = 999.3333
$ws.Range("N82").Value = -1721.3333

$ws.Range("H85").Value = 1189
$ws.Range("J85").Value = 999.3333
$ws.Range("L85").Value = 999.3333
$ws.Range("N85").Value = -3495.3333

$ws.Range("H113").Value = 1343.2778
$ws.Range("I113").Value = 1249.375
$ws.Range("K113").Value = 1249.375
$ws.Range("M113").Value = 920.625

$ws.Range("H122").Value = 9790.647000000001
$ws.Range("I122").Value = 12570.091
$ws.Range("K122").Value = 37710.273
$ws.Range("M122").Value = -35260.273

$ws.Range("H135").Value = 47734.617
$ws.Range("I135").Value = 29687.5
$ws.Range("J135").Value = 105485.4
$ws.Range("K135").Value = 29687.5
$ws.Range("L135").Value = 105485.4
$ws.Range("M135").Value = -24617.5
$ws.Range("N135").Value = -115625.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 11000
$ws.Range("J28").Value = 11000
$ws.Range("L28").Value = 11000
$ws.Range("N28").Value = -11696

$ws.Range("H122").Value = 1901.3
$ws.Range("I122").Value = 2223.6
$ws.Range("K122").Value = 6670.799999999999
$ws.Range("M122").Value = -4220.799999999999
